$wb = $excel.ActiveWorkbook

$wsMS = $wb.Worksheets.Item("MSData")
$wsFlux = $wb.Worksheets.Item("FluxData")
$wsTracer = $wb.Worksheets.Item("TracerData")

# ------------------------------------------------------------------
# FluxData: insert a new fragment row ("EX_glc__D_e.f") right after
# BIOMASS.f (row 2), which pushes the existing rows 3-36 down to 4-37.
# ------------------------------------------------------------------
$wsFlux.Rows.Item(3).Insert()

# Updated error value for the BIOMASS.f row
$wsFlux.Range("C2").Value = 0.0001

# New row 3: EX_glc__D_e.f
$wsFlux.Range("A3").Value = "EX_glc__D_e.f"
$wsFlux.Range("B3").Value = 7.60544986398385
$wsFlux.Range("C3").Value = 3.77142138635765
$wsFlux.Rows.Item(3).RowHeight = 13.8

# Row 4 (previously row 3, EX_c5sugal_e.f) gets new flux/error values
$wsFlux.Range("B4").Value = 0.007523148148148
$wsFlux.Range("C4").Value = 0.016417611948373
$wsFlux.Rows.Item(4).RowHeight = 15

# Former last row (36) keeps its value but grows to the "normal" row height
$wsFlux.Rows.Item(36).RowHeight = 15

# New trailing row 37: DIL_ade_d1.f (default placeholder values)
$wsFlux.Range("A37").Value = "DIL_ade_d1.f"
$wsFlux.Range("B37").Value = 100
$wsFlux.Range("C37").Value = 0.0001
$wsFlux.Rows.Item(37).RowHeight = 13.8

# ------------------------------------------------------------------
# Sheet view / selection / zoom / active-tab bookkeeping
# ------------------------------------------------------------------

# MSData: no longer the selected tab, zoom 55 -> 100, selection shrinks to A2
$wsMS.Activate()
$excel.ActiveWindow.Zoom = 100
$wsMS.Range("A2").Select()

# TracerData: zoom 55 -> 100, selection shrinks to A1
$wsTracer.Activate()
$excel.ActiveWindow.Zoom = 100
$wsTracer.Range("A1").Select()

# FluxData becomes the selected/active tab, zoom 55 -> 100, top-left back to A1,
# selection shrinks to A2
$wsFlux.Activate()
$excel.ActiveWindow.Zoom = 100
$wsFlux.Range("A1").Select()
$excel.ActiveWindow.ScrollRow = 1
$wsFlux.Range("A2").Select()
